$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$cell = $ws.Cells.Item(2, 35)
$v = $cell.Value2
$s = $v.ToString("G17")
Write-Host ("value2 full=" + $s)
